# Re-run SGNN to annotate dialog acts following clean up work to the original transcripts.
# Updates DAMSLTag (column I) and DialogAct (column J) values for the specified rows.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$updates = @(
    @{ Row = 5; I = 'sv'; J = 'Statement-opinion' },
    @{ Row = 21; I = 'sd'; J = 'Statement-non-opinion' },
    @{ Row = 22; I = 'b'; J = 'Acknowledge (Backchannel)' },
    @{ Row = 27; I = 'sd'; J = 'Statement-non-opinion' },
    @{ Row = 29; I = 'sd'; J = 'Statement-non-opinion' },
    @{ Row = 35; I = 'aa'; J = 'Agree/Accept' },
    @{ Row = 43; I = 'aa'; J = 'Agree/Accept' },
    @{ Row = 44; I = 'ba'; J = 'Appreciation' },
    @{ Row = 50; I = '%'; J = 'Uninterpretable' },
    @{ Row = 51; I = '%'; J = 'Uninterpretable' },
    @{ Row = 56; I = 'sd'; J = 'Statement-non-opinion' },
    @{ Row = 57; I = 'b'; J = 'Acknowledge (Backchannel)' },
    @{ Row = 60; I = 'b'; J = 'Acknowledge (Backchannel)' },
    @{ Row = 67; I = 'aa'; J = 'Agree/Accept' },
    @{ Row = 68; I = 'sd'; J = 'Statement-non-opinion' },
    @{ Row = 70; I = 'b'; J = 'Acknowledge (Backchannel)' },
    @{ Row = 76; I = 'b'; J = 'Acknowledge (Backchannel)' },
    @{ Row = 77; I = 'sd'; J = 'Statement-non-opinion' },
    @{ Row = 78; I = 'b'; J = 'Acknowledge (Backchannel)' },
    @{ Row = 95; I = 'sd'; J = 'Statement-non-opinion' },
    @{ Row = 99; I = '%'; J = 'Uninterpretable' },
    @{ Row = 100; I = 'sv'; J = 'Statement-opinion' },
    @{ Row = 103; I = 'sv'; J = 'Statement-opinion' },
    @{ Row = 107; I = 'b'; J = 'Acknowledge (Backchannel)' },
    @{ Row = 119; I = 'b'; J = 'Acknowledge (Backchannel)' },
    @{ Row = 131; I = 'aa'; J = 'Agree/Accept' },
    @{ Row = 142; I = 'b'; J = 'Acknowledge (Backchannel)' },
    @{ Row = 145; I = 'b'; J = 'Acknowledge (Backchannel)' },
    @{ Row = 147; I = 'b'; J = 'Acknowledge (Backchannel)' },
    @{ Row = 166; I = 'b'; J = 'Acknowledge (Backchannel)' },
    @{ Row = 179; I = 'aa'; J = 'Agree/Accept' },
    @{ Row = 197; I = 'aa'; J = 'Agree/Accept' },
    @{ Row = 208; I = 'sd'; J = 'Statement-non-opinion' },
    @{ Row = 216; I = 'sd'; J = 'Statement-non-opinion' },
    @{ Row = 220; I = 'sv'; J = 'Statement-opinion' },
    @{ Row = 225; I = 'b'; J = 'Acknowledge (Backchannel)' },
    @{ Row = 226; I = 'sv'; J = 'Statement-opinion' },
    @{ Row = 232; I = 'b'; J = 'Acknowledge (Backchannel)' },
    @{ Row = 235; I = 'aa'; J = 'Agree/Accept' },
    @{ Row = 293; I = 'aa'; J = 'Agree/Accept' },
    @{ Row = 296; I = 'sd'; J = 'Statement-non-opinion' },
    @{ Row = 329; I = 'b'; J = 'Acknowledge (Backchannel)' },
    @{ Row = 336; I = 'sd'; J = 'Statement-non-opinion' },
    @{ Row = 342; I = 'sd'; J = 'Statement-non-opinion' },
    @{ Row = 351; I = 'sv'; J = 'Statement-opinion' },
    @{ Row = 353; I = 'sd'; J = 'Statement-non-opinion' },
    @{ Row = 372; I = 'b'; J = 'Acknowledge (Backchannel)' }
)

foreach ($u in $updates) {
    $ws.Range("I$($u.Row)").Value = $u.I
    $ws.Range("J$($u.Row)").Value = $u.J
}
